# ---------------------------------------------------------------------------
# Weekly CompStat refresh (66th Precinct): new week, new crime tallies.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: volume/issue number and reporting week -----------------------
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# --- Crime-complaint table (rows 15-28, 33): updated weekly/28-day/YTD counts + %chg ---

$t = $ws.Cells.Item(15,3)
$t.Value = 1
$ws.Range("J14").Copy()
$t.PasteSpecial(-4122)
$ws.Cells.Item(15,6).Value = 3
$ws.Cells.Item(15,9).Value = 11
$ws.Cells.Item(15,11).Value = -8.333333333333
$ws.Cells.Item(15,12).Value = -21.428571428571
$ws.Cells.Item(15,13).Value = 57.142857142857
$ws.Cells.Item(15,14).Value = -21.428571428571
$ws.Cells.Item(16,3).Value = 2
$t = $ws.Cells.Item(16,4)
$t.Value = 1
$ws.Range("J14").Copy()
$t.PasteSpecial(-4122)
$t = $ws.Cells.Item(16,5)
$t.Value = 100
$ws.Range("K14").Copy()
$t.PasteSpecial(-4122)
$ws.Cells.Item(16,6).Value = 4
$ws.Cells.Item(16,7).Value = 7
$ws.Cells.Item(16,8).Value = -42.857142857142
$ws.Cells.Item(16,9).Value = 67
$ws.Cells.Item(16,10).Value = 78
$ws.Cells.Item(16,11).Value = -14.102564102564
$ws.Cells.Item(16,12).Value = 15.517241379310
$ws.Cells.Item(16,13).Value = -44.166666666666
$ws.Cells.Item(16,14).Value = -86.298568507157
$ws.Cells.Item(17,3).Value = 2
$ws.Cells.Item(17,4).Value = 6
$ws.Cells.Item(17,5).Value = -66.666666666666
$ws.Cells.Item(17,6).Value = 19
$ws.Cells.Item(17,7).Value = 17
$ws.Cells.Item(17,8).Value = 11.764705882352
$ws.Cells.Item(17,9).Value = 156
$ws.Cells.Item(17,10).Value = 146
$ws.Cells.Item(17,11).Value = 6.849315068493
$ws.Cells.Item(17,12).Value = 6.122448979591
$ws.Cells.Item(17,13).Value = 41.818181818181
$ws.Cells.Item(17,14).Value = -31.578947368421
$ws.Cells.Item(18,3).Value = 1
$ws.Cells.Item(18,4).Value = 4
$ws.Cells.Item(18,5).Value = -75
$ws.Cells.Item(18,7).Value = 14
$ws.Cells.Item(18,8).Value = -64.285714285714
$ws.Cells.Item(18,9).Value = 85
$ws.Cells.Item(18,10).Value = 82
$ws.Cells.Item(18,11).Value = 3.658536585365
$ws.Cells.Item(18,12).Value = 1.190476190476
$ws.Cells.Item(18,13).Value = -67.181467181467
$ws.Cells.Item(18,14).Value = -92.741246797608
$ws.Cells.Item(19,3).Value = 5
$ws.Cells.Item(19,4).Value = 16
$ws.Cells.Item(19,5).Value = -68.75
$ws.Cells.Item(19,6).Value = 34
$ws.Cells.Item(19,7).Value = 46
$ws.Cells.Item(19,8).Value = -26.086956521739
$ws.Cells.Item(19,9).Value = 301
$ws.Cells.Item(19,10).Value = 375
$ws.Cells.Item(19,11).Value = -19.733333333333
$ws.Cells.Item(19,12).Value = -23.409669211195
$ws.Cells.Item(19,13).Value = 11.481481481481
$ws.Cells.Item(19,14).Value = -32.207207207207
$ws.Cells.Item(20,3).Value = 2
$t = $ws.Cells.Item(20,4)
$t.NumberFormat = "@"
$t.Value = "0"
$ws.Range("C14").Copy()
$t.PasteSpecial(-4122)
$t = $ws.Cells.Item(20,5)
$t.NumberFormat = "@"
$t.Value = "***.*"
$ws.Range("E14").Copy()
$t.PasteSpecial(-4122)
$ws.Cells.Item(20,6).Value = 17
$ws.Cells.Item(20,7).Value = 11
$ws.Cells.Item(20,8).Value = 54.545454545454
$ws.Cells.Item(20,9).Value = 158
$ws.Cells.Item(20,11).Value = 28.455284552845
$ws.Cells.Item(20,12).Value = 50.476190476190
$ws.Cells.Item(20,13).Value = 61.224489795918
$ws.Cells.Item(20,14).Value = -87.227162489894
$ws.Cells.Item(21,3).Value = 13
$ws.Cells.Item(21,4).Value = 27
$ws.Cells.Item(21,5).Value = -51.851851851851
$ws.Cells.Item(21,6).Value = 82
$ws.Cells.Item(21,7).Value = 95
$ws.Cells.Item(21,8).Value = -13.684210526315
$ws.Cells.Item(21,9).Value = 778
$ws.Cells.Item(21,10).Value = 820
$ws.Cells.Item(21,11).Value = -5.121951219512
$ws.Cells.Item(21,12).Value = -3.113325031133
$ws.Cells.Item(21,13).Value = -9.953703703703
$ws.Cells.Item(21,14).Value = -78.334725703146
$t = $ws.Cells.Item(22,3)
$t.Value = 1
$ws.Range("J14").Copy()
$t.PasteSpecial(-4122)
$t = $ws.Cells.Item(22,4)
$t.Value = 1
$ws.Range("J14").Copy()
$t.PasteSpecial(-4122)
$t = $ws.Cells.Item(22,5)
$t.Value = 0
$ws.Range("K14").Copy()
$t.PasteSpecial(-4122)
$ws.Cells.Item(22,6).Value = 2
$t = $ws.Cells.Item(22,7)
$t.Value = 1
$ws.Range("J14").Copy()
$t.PasteSpecial(-4122)
$t = $ws.Cells.Item(22,8)
$t.Value = 100
$ws.Range("K14").Copy()
$t.PasteSpecial(-4122)
$ws.Cells.Item(22,9).Value = 11
$ws.Cells.Item(22,10).Value = 11
$ws.Cells.Item(22,12).Value = 37.5
$ws.Cells.Item(22,13).Value = -31.25
$ws.Cells.Item(24,3).Value = 17
$ws.Cells.Item(24,4).Value = 14
$ws.Cells.Item(24,5).Value = 21.428571428571
$ws.Cells.Item(24,7).Value = 66
$ws.Cells.Item(24,8).Value = 6.060606060606
$ws.Cells.Item(24,9).Value = 671
$ws.Cells.Item(24,10).Value = 669
$ws.Cells.Item(24,11).Value = 0.298953662182
$ws.Cells.Item(24,12).Value = -11.594202898550
$ws.Cells.Item(24,13).Value = 7.36
$ws.Cells.Item(25,4).Value = 5
$ws.Cells.Item(25,5).Value = -40
$ws.Cells.Item(25,7).Value = 20
$ws.Cells.Item(25,8).Value = -55
$ws.Cells.Item(25,9).Value = 104
$ws.Cells.Item(25,10).Value = 169
$ws.Cells.Item(25,11).Value = -38.461538461538
$ws.Cells.Item(25,12).Value = -54.185022026431
$ws.Cells.Item(26,3).Value = 10
$ws.Cells.Item(26,4).Value = 10
$ws.Cells.Item(26,5).Value = 0
$ws.Cells.Item(26,6).Value = 30
$ws.Cells.Item(26,8).Value = -14.285714285714
$ws.Cells.Item(26,9).Value = 260
$ws.Cells.Item(26,10).Value = 289
$ws.Cells.Item(26,11).Value = -10.034602076124
$ws.Cells.Item(26,12).Value = 8.786610878661
$ws.Cells.Item(26,13).Value = -10.958904109589
$t = $ws.Cells.Item(27,3)
$t.Value = 1
$ws.Range("J14").Copy()
$t.PasteSpecial(-4122)
$ws.Cells.Item(27,6).Value = 3
$t = $ws.Cells.Item(27,7)
$t.NumberFormat = "@"
$t.Value = "0"
$ws.Range("C14").Copy()
$t.PasteSpecial(-4122)
$t = $ws.Cells.Item(27,8)
$t.NumberFormat = "@"
$t.Value = "***.*"
$ws.Range("E14").Copy()
$t.PasteSpecial(-4122)
$ws.Cells.Item(27,9).Value = 13
$ws.Cells.Item(27,11).Value = -40.909090909090
$ws.Cells.Item(27,12).Value = -18.75
$ws.Cells.Item(28,3).Value = 2
$t = $ws.Cells.Item(28,4)
$t.NumberFormat = "@"
$t.Value = "0"
$ws.Range("C14").Copy()
$t.PasteSpecial(-4122)
$t = $ws.Cells.Item(28,5)
$t.NumberFormat = "@"
$t.Value = "***.*"
$ws.Range("E14").Copy()
$t.PasteSpecial(-4122)
$ws.Cells.Item(28,6).Value = 5
$ws.Cells.Item(28,8).Value = -16.666666666666
$ws.Cells.Item(28,9).Value = 40
$ws.Cells.Item(28,11).Value = -13.043478260869
$ws.Cells.Item(28,12).Value = -18.367346938775
$t = $ws.Cells.Item(33,6)
$t.NumberFormat = "@"
$t.Value = "0"
$ws.Range("C14").Copy()
$t.PasteSpecial(-4122)
$ws.Cells.Item(33,8).Value = -100

# --- Column widths: Week-to-Date %Chg columns (E, H) widened to fit new values ---
$ws.Columns.Item(5).ColumnWidth = 6.7
$ws.Columns.Item(8).ColumnWidth = 6.7
